$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The status text "Ready for handoff" changed to "Handback transform failed"
# everywhere it is used (shared string reused across sheets).
$oldStatus = "Ready for handoff"
$newStatus = "Handback transform failed"

if ($wsOverview.Range("E3").Value2 -eq $oldStatus) { $wsOverview.Range("E3").Value = $newStatus }
if ($wsOverview.Range("F3").Value2 -eq $oldStatus) { $wsOverview.Range("F3").Value = $newStatus }
if ($wsZhCn.Range("C3").Value2 -eq $oldStatus) { $wsZhCn.Range("C3").Value = $newStatus }
if ($wsDeDe.Range("C3").Value2 -eq $oldStatus) { $wsDeDe.Range("C3").Value = $newStatus }

# Populate the "Error Detail" column (P) for row 3 on the zh-cn and de-de
# sheets with the handback-transform error message, reporting the
# handback/handoff file-name mismatch for each locale.
$wsZhCn.Range("P3").Value = "Handback file name: y1gf3ds0.qio is different with handoff file name: e1952233-b6a0-4378-8b71-a236379f8f50.0b125d5997fe0224b631b6161a15fe021aaaaa60.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: y1gf3ds0.qio is different with handoff file name: e1952233-b6a0-4378-8b71-a236379f8f50.0b125d5997fe0224b631b6161a15fe021aaaaa60.de-de."

# Widen column P (Error Detail) on both locale sheets so the new, longer
# error text is readable - match the width already used by other wide
# columns in the sheet (e.g. column A, which is stored as width 40).
$wideWidth = $wsZhCn.Columns(1).ColumnWidth
$wsZhCn.Columns(16).ColumnWidth = $wideWidth
$wsDeDe.Columns(16).ColumnWidth = $wideWidth
